$d = $word.ActiveDocument

# --- locate the target field (the "m:'...'.asParagraph().setAlignment('NOT_EXISTING')" field) ---
$f = $d.Fields.Item(1)
$code = $f.Code.Text                       # e.g. " m:'some text'.asParagraph().setAlignment('NOT_EXISTING') "
$inner = $code.Trim()                      # "m:'some text'.asParagraph().setAlignment('NOT_EXISTING')"

# the paragraph that contains the field
$fieldParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Start -le $f.Code.Start -and $d.Paragraphs($i).Range.End -ge $f.Code.End) {
        $fieldParaIndex = $i
    }
}

# delete the field entirely (removes fldChar begin/end, all instrText runs and the _GoBack bookmark)
$f.Delete()

$para = $d.Paragraphs($fieldParaIndex).Range
$paraStart = $para.Start

# literal text that replaces the field: "{" + field code (trimmed) + "}"
$full = "{" + $inner + "}"
$para2 = $d.Paragraphs($fieldParaIndex).Range
$para2.Text = $full

# give the whole run the same language formatting the original instrText runs carried
$whole = $d.Range($paraStart, $paraStart + $full.Length)
$whole.LanguageID = "en-US"

# --- split the single merged run back into the original per-token runs ---
# list of (length) for every token, in order; a short Bold-toggle forces Word to keep
# each token in its own <w:r> without altering the (already-applied) language formatting.
$lengths = @(1,1,2,9,1,14,1,3,9,1,1,12,1,1,1)

$pos = $paraStart
foreach ($len in $lengths) {
    $r = $d.Range($pos, $pos + $len)
    $r.Bold = 1
    $r.Bold = 0
    $pos = $pos + $len
}

# the very last token ("}") must not carry any rPr at all (matches the original
# un-formatted trailing run) -- clear the language we set on it earlier.
$lastStart = $paraStart + $full.Length - 1
$lastRun = $d.Range($lastStart, $paraStart + $full.Length)
$lastRun.LanguageID = 0

# --- re-insert the hidden _GoBack bookmark exactly where it used to sit: right
# after "NOT_EXISTING" and before the following "'" ---
$bmPos = $paraStart + 1 + 1 + 2 + 9 + 1 + 14 + 1 + 3 + 9 + 1 + 1 + 12
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))
